$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")

$ws.Range("J2").Value = 416
$ws.Range("K2").Value = 624
$ws.Range("J3").Value = 352
$ws.Range("K3").Value = 320
$ws.Range("K4").Value = 304
$ws.Range("J5").Value = 736
$ws.Range("K5").Value = 1232
$ws.Range("J6").Value = 160
$ws.Range("K6").Value = 688
$ws.Range("J7").Value = 352
$ws.Range("K7").Value = 672
$ws.Range("J8").Value = 672
$ws.Range("K8").Value = 64
$ws.Range("J9").Value = 224
$ws.Range("K9").Value = 464
$ws.Range("J18").Value = 576
$ws.Range("K18").Value = 672
$ws.Range("J19").Value = 1216
$ws.Range("K19").Value = 448
$ws.Range("J20").Value = 624
$ws.Range("K20").Value = 80
$ws.Range("J21").Value = 480
$ws.Range("K21").Value = 800
$ws.Range("J22").Value = 912
$ws.Range("K22").Value = 352
$ws.Range("J23").Value = 592
$ws.Range("K23").Value = 208
$ws.Range("J24").Value = 1072
$ws.Range("K24").Value = 1104
$ws.Range("J25").Value = 720
$ws.Range("K25").Value = 224
$ws.Range("J26").Value = 224
$ws.Range("K26").Value = 1088
$ws.Range("J27").Value = 288
$ws.Range("K27").Value = 192
$ws.Range("J28").Value = 880
$ws.Range("K28").Value = 1104
$ws.Range("J29").Value = 64
$ws.Range("K29").Value = 352
$ws.Range("J30").Value = 448
$ws.Range("K30").Value = 528
$ws.Range("J31").Value = 704
$ws.Range("K31").Value = 688
$ws.Range("J32").Value = 816
$ws.Range("K32").Value = 1120
$ws.Range("J33").Value = 800
$ws.Range("J34").Value = 80
$ws.Range("K34").Value = 656
$ws.Range("J36").Value = 576
$ws.Range("K36").Value = 336
$ws.Range("J37").Value = 816
$ws.Range("K37").Value = 608
$ws.Range("J38").Value = 400
$ws.Range("K38").Value = 608
$ws.Range("J39").Value = 992
$ws.Range("K39").Value = 1136
$ws.Range("J40").Value = 624
$ws.Range("K40").Value = 128
$ws.Range("J41").Value = 944
$ws.Range("K41").Value = 368
$ws.Range("J42").Value = 944
$ws.Range("K42").Value = 464
$ws.Range("J43").Value = 128
$ws.Range("K43").Value = 784
$ws.Range("J44").Value = 496
$ws.Range("K44").Value = 928
$ws.Range("J45").Value = 560
$ws.Range("K45").Value = 832
$ws.Range("J46").Value = 128
$ws.Range("K46").Value = 128
$ws.Range("J47").Value = 224
$ws.Range("K47").Value = 128
$ws.Range("J48").Value = 672
$ws.Range("K48").Value = 480
$ws.Range("J49").Value = 1168
$ws.Range("K49").Value = 608
$ws.Range("J50").Value = 688
$ws.Range("K50").Value = 432
$ws.Range("J51").Value = 176
$ws.Range("K51").Value = 560
$ws.Range("J52").Value = 1152
$ws.Range("K52").Value = 240
$ws.Range("J53").Value = 464
$ws.Range("K53").Value = 256
$ws.Range("J54").Value = 624
$ws.Range("K54").Value = 688
$ws.Range("J55").Value = 64
$ws.Range("K55").Value = 496
$ws.Range("J56").Value = 960
$ws.Range("K56").Value = 64
